$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.834.64'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '1.638.83'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''214.43'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('D9').Value = '''0.0612'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').Value = '''19.44'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').Value = '1.867.60'
$ws.Range('E12').Value = '  +2.14%  '
$ws.Range('D13').Value = '1.634.74'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').Value = '''0.516'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').Value = '''64.57'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '''241.72'
$ws.Range('E17').Value = '  +6.40%  '
$ws.Range('D18').Value = '26.802.90'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').Value = '''7.87'
$ws.Range('E19').Value = '  +3.10%  '
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').Value = '''4.39'
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('D23').Value = '''2.24'
$ws.Range('E23').Value = '  +3.34%  '
$ws.Range('E24').Value = '  +3.08%  '
$ws.Range('D25').Value = '''145.94'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').Value = '''0.0495'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = '''1.17'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''3.28'
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.510.36'
$ws.Range('E33').Value = '  +4.62%  '
$ws.Range('E34').Value = '  +2.31%  '
$ws.Range('E35').Value = '  +6.50%  '
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('D37').Value = '''0.574'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('D39').Value = '''0.859'
$ws.Range('E39').Value = '  +4.17%  '
$ws.Range('D40').Value = '''5.96'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('D43').Value = '''64.17'
$ws.Range('E43').Value = '  +5.46%  '
$ws.Range('D44').Value = '1.778.26'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('D47').Value = '''90.39'
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('E48').Value = '  +3.02%  '
$ws.Range('D49').Value = '''0.0977'
$ws.Range('E49').Value = '  +2.74%  '
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('D51').Value = '''7.51'
$ws.Range('E51').Value = '  +1.68%  '
